# Rename the "temperature_c" sheet to "temperature" and make it the
# active/selected tab (moving the active tab away from "genotype").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("temperature_c")
$ws.Name = "temperature"
$ws.Activate()
